# optimise the procedure of change job
#
# The quest description text in column I ("Descript") used a "G|<name>||..."
# marker to reference the quest-giving NPC. This is changed to the more
# explicit "#npc|<name>||..." marker (matching the |#scene|, |#event|,
# |#item| markers already used elsewhere in the same strings), for every
# quest row (rows 4-15) on the "Quest" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$descriptions = @{
    4  = "#npc|乔斯基||让你去找|#npc|瓦里斯||，并向他学习战斗的技巧，提升自己。"
    5  = "#npc|乔斯基||希望你可以去到|#scene|布萨特高塔||，并在那里历练自己。"
    6  = "#npc|塞尼斯||希望你帮助他，进入附近的|#event|狼穴||并消灭狼群，如果失败了，你可以尝试反复进出本地图重试。"
    7  = "#npc|塞巴斯恰恩||的宠物|#npc|恰恰||不见了，他希望你可以帮助他找回宠物。"
    8  = "#npc|奥莱伊李||觉得目前的田野很不安全，希望你可以帮助他清理当前场景的所有特殊事件。"
    9  = "#npc|玛莎||告诉你，他的孩子，前几天走失了。如果你碰到了这个孩子，一定要记得把他带回来"
    10 = "#npc|科迪||告诉你，附近的|#scene|村落遗迹||中，隐藏着一些秘密，你可以去调查下，说不定会有意外的收获。"
    11 = "#npc|乔斯基||需要一些罂粟花苗，让你寻找收集一些。应该可以从|#event|罂粟花田||中找到。"
    12 = "#npc|塞尼斯||告诉你，再附近的森林深处，有一只神兽|#event|穷奇||，如果你可以找到并击败他，会得到丰厚的回报。"
    13 = "#npc|奥莱伊李||教了你种植植物的方法，||你需要到附近的田地里，种下|#item|豌豆种子||。并在收获后，把果实带交给他。"
    14 = "#npc|塞巴斯恰恩||委托你找到可恶的|#npc|科迪||，并通过卡牌战斗的方式击败他。|#npc|科迪||经常欺负周边的村民，非常可恶。"
    15 = "#npc|科迪||希望你可以找到并战胜|#scene|村落遗迹||中的|#event|蛮王之灵||。"
}

foreach ($row in $descriptions.Keys) {
    $ws.Range("I$row").Value2 = $descriptions[$row]
}

# Update the view state left over from editing: scrolled so column I is the
# left-most visible column, with the final quest row's description selected.
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("I1").Column
$ws.Range("I8").Select()
